$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (order matches shared-string insertion order: A, D, B, C, F)
$ws.Range("A5").Value = "Frozen BERT, fine tuned classification head"
$ws.Range("D5").Value = "~2 iter/sec, 375 iters"
$ws.Range("B5").Value = "2.5 minutes"
$ws.Range("C5").Value = 0.2
$ws.Range("F5").Value = "needs to relearn embeddings for entity labels?"

# Row 6 (order: A, B, C, D)
$ws.Range("A6").Value = "Fine tuned embeddings and classification"
$ws.Range("B6").Value = "4 minutes"
$ws.Range("C6").Value = 0.2
$ws.Range("D6").Value = "~1.5 iter/sec, 375 iters"

# Update selection like the diff shows (activeCell D18)
$ws.Range("D18").Select()
